$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Belief about GCS"
$ws.Range("B2").Value = 60.941935483871
$ws.Range("C2").Value = 55.1875825627477
$ws.Range("D2").Value = 63.3075506445672
$ws.Range("E2").Value = 57.4937888198758
$ws.Range("F2").Value = 50.9511918274688

# Row 3 - "Belief about NR"
$ws.Range("B3").Value = 60.3854838709677
$ws.Range("C3").Value = 52.672391017173
$ws.Range("D3").Value = 62.6279926335175
$ws.Range("E3").Value = 59.3400621118012
$ws.Range("F3").Value = 54.1278849791903
